$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @("永辉超市", "C摩尔-U", "万科A")
    3  = @("万  科Ａ", "永辉超市", "永辉超市")
    4  = @("航天发展", "万  科Ａ", "航天发展")
    5  = @("C摩尔-U", "中科曙光", "龙洲股份")
    6  = @("东百集团", "龙洲股份", "中科曙光")
    7  = @("中科曙光", "航天发展", "海马汽车")
    8  = @("龙洲股份", "东百集团", "摩尔线程")
    9  = @("海马汽车", "海马汽车", "平潭发展")
    10 = @("西部材料", "达华智能", "特发信息")
    11 = @("安记食品", "天通股份", "安记食品")
    12 = @("平潭发展", "合力泰", "海王生物")
    13 = @("海王生物", "平潭发展", "东百集团")
    14 = @("海光信息", "海王生物", "海南发展")
    15 = @("天通股份", "海南发展", "安妮股份")
    16 = @("海南发展", "安记食品", "博纳影业")
    17 = @("达华智能", "华夏幸福", "西部材料")
    18 = @("雷科防务", "神农种业", "厦门港务")
    19 = @("神农种业", "雷科防务", "海欣食品")
    20 = @("罗 牛 山", "罗 牛 山", "实达集团")
    21 = @("特发信息", "特发信息", "合富中国")
}

foreach ($row in $values.Keys | Sort-Object) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
}
